$d = $word.ActiveDocument

$replacements = @(
    @("313÷4=78, 1", "416÷2=208, 0"),
    @("443÷4=110, 3", "445÷3=148, 1"),
    @("249÷2=124, 1", "186÷5=37, 1"),
    @("496÷5=99, 1", "287÷2=143, 1"),
    @("162÷3=54, 0", "360÷3=120, 0"),
    @("881÷7=125, 6", "437÷3=145, 2"),
    @("704÷3=234, 2", "599÷2=299, 1"),
    @("804÷4=201, 0", "964÷4=241, 0"),
    @("794÷7=113, 3", "554÷8=69, 2"),
    @("323÷7=46, 1", "831÷3=277, 0"),
    @("861÷6=143, 3", "671÷2=335, 1"),
    @("506÷9=56, 2", "705÷4=176, 1"),
    @("525÷6=87, 3", "188÷2=94, 0"),
    @("249÷8=31, 1", "431÷3=143, 2"),
    @("129÷7=18, 3", "935÷7=133, 4"),
    @("698÷2=349, 0", "334÷8=41, 6"),
    @("629÷6=104, 5", "419÷4=104, 3"),
    @("955÷5=191, 0", "662÷9=73, 5"),
    @("917÷3=305, 2", "720÷3=240, 0"),
    @("736÷2=368, 0", "147÷6=24, 3"),
    @("198÷8=24, 6", "321÷3=107, 0"),
    @("259÷3=86, 1", "164÷3=54, 2"),
    @("386÷7=55, 1", "418÷3=139, 1"),
    @("127÷7=18, 1", "568÷9=63, 1"),
    @("992÷6=165, 2", "717÷4=179, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
